$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("My Data")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet1 ("My Data") ---
# Header row now at row 1 (was row 4)
$headers = @("Code","BusName","BusPhone","ContactFirst","ContactLast","Additional")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws1.Cells.Item(1, $c + 1).Value = $headers[$c]
}

# Data rows 2..11 for codes 123..132
for ($i = 0; $i -lt 10; $i++) {
    $code = 123 + $i
    $row = 2 + $i
    $ws1.Cells.Item($row, 1).Value = $code
    $ws1.Cells.Item($row, 2).Value = "$code Name"
    $ws1.Cells.Item($row, 3).Value = 1230897
    $ws1.Cells.Item($row, 4).Value = "John"
    $ws1.Cells.Item($row, 5).Value = "Citizen"
    $ws1.Cells.Item($row, 6).Value = "Testing Testing Testing"
}

$ws1.Range("A1:F11").Select()
try { $ws1.Range("A1:F11").SetPhonetic() } catch {}

# --- Sheet2 ---
$ws2.Cells.Item(3, 3).Value = "as"
$ws2.Cells.Item(1, 2).Value = "a"
$ws2.Cells.Item(2, 2).Value = "a"
$ws2.Cells.Item(3, 2).Value = "a"
$ws2.Cells.Item(5, 2).Value = "as"
$ws2.Cells.Item(8, 4).Value = "d"

for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws2.Cells.Item(9, $c + 1).Value = $headers[$c]
}

for ($i = 0; $i -lt 3; $i++) {
    $code = 123 + $i
    $row = 10 + $i
    $ws2.Cells.Item($row, 1).Value = $code
    $ws2.Cells.Item($row, 2).Value = "$code Name"
    $ws2.Cells.Item($row, 3).Value = 1230897
    $ws2.Cells.Item($row, 4).Value = "John"
    $ws2.Cells.Item($row, 5).Value = "Citizen"
    $ws2.Cells.Item($row, 6).Value = "Testing Testing Testing"
}

$ws2.Cells.Item(17, 7).Select()
$ws2.Activate()
